$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.892.15'
$ws.Range("E2").Value = '  -1.65%  '
$ws.Range("D3").Value = '2.302.36'
$ws.Range("E3").Value = '  -1.62%  '
$ws.Range("E4").Value = '  -0.05%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.60'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  -1.45%  '
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.03'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  -5.03%  '
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.505'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  -1.40%  '
$ws.Range("E8").Value = '  +0.04%  '
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.495'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  -4.15%  '
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.58'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  -5.00%  '
$ws.Range("E11").Value = '  -0.49%  '
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '49.21'
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = '  -5.01%  '
$ws.Range("E13").Value = '  +2.13%  '
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.87'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  +8.43%  '
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.78'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  -0.75%  '
$ws.Range("D16").Value = '2.661.71'
$ws.Range("E16").Value = '  -1.79%  '
$ws.Range("D17").Value = '2.313.98'
$ws.Range("E17").Value = '  -1.46%  '
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.807'
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  -0.07%  '
$ws.Range("D19").Value = '42.851.42'
$ws.Range("E19").Value = '  -1.54%  '
$ws.Range("D20").Value = '0.0₃0900'
$ws.Range("E20").Value = '  -0.98%  '
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.58'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  -2.28%  '
$ws.Range("E22").Value = '  -2.19%  '
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.24'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  -1.38%  '
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '236.09'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  -1.00%  '
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.01'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("E26").Value = '  +0.05%  '
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.45'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  -3.90%  '
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.66'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  -1.49%  '
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.28'
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = '  +4.74%  '
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.46'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  +0.17%  '
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.82'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  -2.48%  '
$ws.Range("E32").Value = '  -1.70%  '
$ws.Range("E33").Value = '  -0.10%  '
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.77'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  +6.05%  '
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.95'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  -2.44%  '
$ws.Range("E36").Value = '  -1.27%  '
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.81'
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  -0.84%  '
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0694'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  -1.96%  '
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.81'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  -3.49%  '
$ws.Range("E40").Value = '  -2.14%  '
$ws.Range("E41").Value = '  -4.56%  '
$ws.Range("E42").Value = '  -2.25%  '
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.34'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  -2.32%  '
$ws.Range("D44").Value = '1.976.23'
$ws.Range("E44").Value = '  -0.80%  '
$ws.Range("E45").Value = '  -1.98%  '
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.65'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  -4.94%  '
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.82'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  -1.56%  '
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.84'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  -4.06%  '
$ws.Range("D49").Value = '2.527.66'
$ws.Range("E49").Value = '  -1.53%  '
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.76'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  -6.95%  '
$ws.Range("E51").Value = '  -6.58%  '
